# 12.4.2.xlsx — add the 2023 (column M) data series, bold the "per-capita"
# summary row, tidy up row/column sizing and wrap the long "share buried"
# label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row gets a bit taller now that the label wraps across 13 cols ---
$ws.Rows.Item(1).RowHeight = 57

# --- Columns A:C share one uniform width instead of three bespoke ones ---
$ws.Columns.Item(1).ColumnWidth = 37.83
$ws.Columns.Item(2).ColumnWidth = 37.83
$ws.Columns.Item(3).ColumnWidth = 37.83

# --- Extend the year header line with 2023 ---
$ws.Range("L3:L8").Copy()
$ws.Range("M3:M8").PasteSpecial(-4122)

$ws.Range("M3").Value = 2023
$ws.Range("M4").Value = 1963.9481143272037
$ws.Range("M5").Value = 14065.6
$ws.Range("M6").Value = 7161.9
$ws.Range("M7").Value = 46.213456944602434
$ws.Range("M8").Value = 0.044790126265498803

# --- Make the "per capita" row (row 4) stand out in bold, like the header ---
$ws.Rows.Item(4).Font.Bold = $true

# --- Wrap the long "share of buried hazardous waste" label and let the row grow ---
$ws.Range("A7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 24

# --- Clear the leftover cell selection from the previous edit session ---
$ws.Range("A1").Select()
